# Updates cryptos list values (Price and Volume(1h) columns) per commit
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'25.989.73"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +0.65%  "
$ws.Range("D3").Value = "'1.739.06"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.01%  "
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'239.89"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.59%  "
$ws.Range("E6").Value = "  +0.13%  "
$ws.Range("D7").Value = "'0.5282"
$ws.Range("D7").ClearFormats()
$ws.Range("D8").Value = "'0.2741"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  -2.27%  "
$ws.Range("D9").Value = "'0.06163"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +0.78%  "
$ws.Range("D10").Value = "'1.740.10"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -0.04%  "
$ws.Range("D11").Value = "'0.07186"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +2.15%  "
$ws.Range("D12").Value = "'15.10"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  -1.27%  "
$ws.Range("D13").Value = "'0.6413"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  -0.95%  "
$ws.Range("D14").Value = "'4.605"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.75%  "
$ws.Range("D15").Value = "'77.52"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +0.72%  "
$ws.Range("E16").Value = "  +0.12%  "
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Value = "'26.015.05"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.76%  "
$ws.Range("D19").Value = "'11.78"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +2.59%  "
$ws.Range("D20").Value = "'0.000006755"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +2.34%  "
$ws.Range("D21").Value = "'1.965.62"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  -0.50%  "
$ws.Range("D22").Value = "'4.359"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +5.33%  "
$ws.Range("D23").Value = "'8.625"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -0.42%  "
$ws.Range("D24").Value = "'5.248"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +2.10%  "
$ws.Range("D25").Value = "'140.30"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.72%  "
$ws.Range("E26").Value = "  +0.29%  "
$ws.Range("D27").Value = "'15.22"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.05%  "
$ws.Range("D28").Value = "'1.767"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -2.63%  "
$ws.Range("D29").Value = "'105.57"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.32%  "
$ws.Range("D30").Value = "'0.08398"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.99%  "
$ws.Range("E31").Value = "  +4.23%  "
$ws.Range("D32").Value = "'3.648"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +6.39%  "
$ws.Range("D33").Value = "'0.04587"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.07%  "
$ws.Range("D34").Value = "'2.652"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +1.63%  "
$ws.Range("D35").Value = "'0.9917"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.63%  "
$ws.Range("D36").Value = "'0.6241"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +1.51%  "
$ws.Range("D37").Value = "'2.702"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.67%  "
$ws.Range("D38").Value = "'0.01603"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +1.17%  "
$ws.Range("D39").Value = "'1.931"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -0.35%  "
$ws.Range("E40").Value = "  +0.13%  "
$ws.Range("D41").Value = "'98.83"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.90%  "
$ws.Range("D42").Value = "'0.3886"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +1.01%  "
$ws.Range("D43").Value = "'0.7496"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +3.02%  "
$ws.Range("D44").Value = "'4.946"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -0.56%  "
$ws.Range("D45").Value = "'0.1143"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +2.00%  "
$ws.Range("E46").Value = "  -1.82%  "
$ws.Range("D47").Value = "'6.214"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -1.04%  "
$ws.Range("D48").Value = "'54.73"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +2.97%  "
$ws.Range("D49").Value = "'30.78"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +2.77%  "
$ws.Range("D50").Value = "'0.3443"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +1.25%  "
$ws.Range("D51").Value = "'7.522"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  -2.35%  "
